$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 51: coin changed from SynthetixNetwork to HuobiToken
$ws.Range("B51").Value = "HuobiToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.70"
$ws.Range("E51").Value = "  +0.18%  "

# Row 2
$ws.Range("D2").Value = "42.197.24"
$ws.Range("E2").Value = "  -0.98%  "

# Row 3
$ws.Range("D3").Value = "2.243.07"
$ws.Range("E3").Value = "  -1.04%  "

# Row 4
$ws.Range("E4").Value = "  +0.27%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.42"
$ws.Range("E5").Value = "  -1.48%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.620"
$ws.Range("E6").Value = "  -3.14%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "74.24"
$ws.Range("E7").Value = "  -2.73%  "

# Row 8
$ws.Range("E8").Value = "  +0.07%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.615"
$ws.Range("E9").Value = "  -4.71%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.11"
$ws.Range("E10").Value = "  +4.55%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0945"
$ws.Range("E11").Value = "  -2.63%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.15"
$ws.Range("E12").Value = "  -2.70%  "

# Row 13
$ws.Range("E13").Value = "  -2.72%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.53"
$ws.Range("E14").Value = "  -3.14%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.854"
$ws.Range("E15").Value = "  -1.30%  "

# Row 16
$ws.Range("D16").Value = "2.274.44"
$ws.Range("E16").Value = "  +0.19%  "

# Row 17
$ws.Range("D17").Value = "42.102.03"
$ws.Range("E17").Value = "  -1.00%  "

# Row 18
$ws.Range("D18").Value = "0.0₃0991"
$ws.Range("E18").Value = "  -0.26%  "

# Row 19
$ws.Range("E19").Value = "  -0.57%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.94"
$ws.Range("E20").Value = "  -0.19%  "

# Row 21
$ws.Range("E21").Value = "  +4.74%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "231.91"
$ws.Range("E22").Value = "  -0.90%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.57"
$ws.Range("E23").Value = "  +30.50%  "

# Row 24
$ws.Range("E24").Value = "  +0.04%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.26"
$ws.Range("E25").Value = "  -0.30%  "

# Row 27
$ws.Range("E27").Value = "  -2.96%  "

# Row 28
$ws.Range("E28").Value = "  -2.27%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "169.01"
$ws.Range("E29").Value = "  +0.84%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.62"
$ws.Range("E30").Value = "  -1.36%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0815"
$ws.Range("E31").Value = "  -5.17%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.119"
$ws.Range("E32").Value = "  -4.05%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "30.26"
$ws.Range("E33").Value = "  -4.19%  "

# Row 34
$ws.Range("E34").Value = "  -1.65%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.15"
$ws.Range("E35").Value = "  +8.42%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.51"
$ws.Range("E36").Value = "  -0.33%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0303"
$ws.Range("E37").Value = "  -1.14%  "

# Row 38
$ws.Range("E38").Value = "  -0.92%  "

# Row 39
$ws.Range("E39").Value = "  -3.21%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.80"
$ws.Range("E40").Value = "  -1.05%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "61.98"
$ws.Range("E41").Value = "  +0.51%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.201"
$ws.Range("E42").Value = "  -2.74%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "107.62"
$ws.Range("E43").Value = "  +0.85%  "

# Row 44
$ws.Range("E44").Value = "  +1.80%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.64"
$ws.Range("E45").Value = "  -2.37%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.996"
$ws.Range("E46").Value = "  -0.31%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.37"
$ws.Range("E47").Value = "  -7.57%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.11"
$ws.Range("E48").Value = "  -3.45%  "

# Row 49
$ws.Range("E49").Value = "  -0.73%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.28"
$ws.Range("E50").Value = "  +1.40%  "
